# Daily attendance processing - 2025-10-31 17:45:23
# Reorders the "Recorded By" (column G) values for the rows listed below so
# that the author names/emails appear in their updated order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Map of new cell value -> list of row numbers that should receive it.
$updates = @{
    "System, system, backup@backdoor.com" = @(2, 29, 56)
    "System, backup@backdoor.com"         = @(4, 5, 8, 31, 32, 35, 58, 59, 62, 83, 84, 85, 109, 110, 111, 135, 136, 137)
    "System, dnasr281@gmail.com"          = @(11, 17, 38, 44, 65, 71, 96, 97, 99, 122, 123, 125, 148, 149, 151)
    "admin@admin.com, dnasr281@gmail.com" = @(90, 116, 142)
}

foreach ($newValue in $updates.Keys) {
    foreach ($row in $updates[$newValue]) {
        $ws.Range("G$row").Value = $newValue
    }
}
